$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "67.393.42"
$ws.Cells.Item(2, 5).Value = "  +1.73%  "

$ws.Cells.Item(3, 4).Value = "3.786.59"
$ws.Cells.Item(3, 5).Value = "  +6.42%  "

$c = $ws.Cells.Item(4, 4)
$c.Value = "'1.01"
$c.Style = "Normal"
$ws.Cells.Item(4, 5).Value = "  +0.76%  "

$c = $ws.Cells.Item(5, 4)
$c.Value = "'415.87"
$c.Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  -0.95%  "

$c = $ws.Cells.Item(6, 4)
$c.Value = "'137.87"
$c.Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +4.32%  "

$ws.Cells.Item(7, 4).Value = "3.776.63"
$ws.Cells.Item(7, 5).Value = "  +6.39%  "

$c = $ws.Cells.Item(8, 4)
$c.Value = "'0.640"
$c.Style = "Normal"
$ws.Cells.Item(8, 5).Value = "  -3.52%  "

$ws.Cells.Item(9, 5).Value = "  +0.08%  "

$c = $ws.Cells.Item(10, 4)
$c.Value = "'0.763"
$c.Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  -2.67%  "

$c = $ws.Cells.Item(11, 4)
$c.Value = "'0.177"
$c.Style = "Normal"
$ws.Cells.Item(11, 5).Value = "  +5.95%  "

$c = $ws.Cells.Item(12, 4)
$c.Value = "'0.0000371"
$c.Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  +30.36%  "

$c = $ws.Cells.Item(13, 4)
$c.Value = "'42.48"
$c.Style = "Normal"
$ws.Cells.Item(13, 5).Value = "  -2.23%  "

$c = $ws.Cells.Item(14, 4)
$c.Value = "'10.24"
$c.Style = "Normal"
$ws.Cells.Item(14, 5).Value = "  +1.28%  "

$ws.Cells.Item(15, 4).Value = "4.405.69"
$ws.Cells.Item(15, 5).Value = "  +6.94%  "

$ws.Cells.Item(16, 5).Value = "  -0.59%  "

$ws.Cells.Item(17, 4).Value = "3.775.41"
$ws.Cells.Item(17, 5).Value = "  +5.94%  "

$c = $ws.Cells.Item(18, 4)
$c.Value = "'20.50"
$c.Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  -0.43%  "

$c = $ws.Cells.Item(19, 4)
$c.Value = "'13.30"
$c.Style = "Normal"
$ws.Cells.Item(19, 5).Value = "  +4.12%  "

$c = $ws.Cells.Item(20, 4)
$c.Value = "'1.12"
$c.Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +0.04%  "

$ws.Cells.Item(21, 4).Value = "67.705.85"
$ws.Cells.Item(21, 5).Value = "  +2.31%  "

$c = $ws.Cells.Item(22, 4)
$c.Value = "'435.78"
$c.Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -3.28%  "

$c = $ws.Cells.Item(23, 4)
$c.Value = "'15.05"
$c.Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +13.83%  "

$c = $ws.Cells.Item(24, 4)
$c.Value = "'89.16"
$c.Style = "Normal"
$ws.Cells.Item(24, 5).Value = "  -1.55%  "

$c = $ws.Cells.Item(25, 4)
$c.Value = "'3.06"
$c.Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  -5.98%  "

$c = $ws.Cells.Item(26, 4)
$c.Value = "'37.39"
$c.Style = "Normal"
$ws.Cells.Item(26, 5).Value = "  +8.99%  "

$c = $ws.Cells.Item(27, 4)
$c.Value = "'3.28"
$c.Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  -3.23%  "

$c = $ws.Cells.Item(28, 4)
$c.Value = "'9.74"
$c.Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  -3.69%  "

$c = $ws.Cells.Item(29, 4)
$c.Value = "'5.17"
$c.Style = "Normal"
$ws.Cells.Item(29, 5).Value = "  +6.89%  "

$c = $ws.Cells.Item(30, 4)
$c.Value = "'12.55"
$c.Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +0.19%  "

$c = $ws.Cells.Item(31, 4)
$c.Value = "'0.122"
$c.Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +3.49%  "

$c = $ws.Cells.Item(32, 4)
$c.Value = "'2.74"
$c.Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  -1.91%  "

$c = $ws.Cells.Item(33, 4)
$c.Value = "'7.20"
$c.Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -1.46%  "

$c = $ws.Cells.Item(34, 4)
$c.Value = "'41.48"
$c.Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +6.42%  "

$c = $ws.Cells.Item(35, 4)
$c.Value = "'0.161"
$c.Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  -1.01%  "

$c = $ws.Cells.Item(36, 4)
$c.Value = "'57.73"
$c.Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +0.13%  "

$c = $ws.Cells.Item(37, 4)
$c.Value = "'0.999"
$c.Style = "Normal"
$ws.Cells.Item(37, 5).Value = "  +0.03%  "

$c = $ws.Cells.Item(38, 4)
$c.Value = "'0.0484"
$c.Style = "Normal"
$ws.Cells.Item(38, 5).Value = "  -4.66%  "

$c = $ws.Cells.Item(39, 4)
$c.Value = "'2.98"
$c.Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  +27.55%  "

$c = $ws.Cells.Item(40, 4)
$c.Value = "'0.145"
$c.Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -2.62%  "

$ws.Cells.Item(41, 2).Value = "FirstDigitalUSD"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$c = $ws.Cells.Item(41, 4)
$c.Value = "'1.00"
$c.Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +0.18%  "

$ws.Cells.Item(42, 2).Value = "PEPE"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(42, 4).Value = "0.0₃0668"
$ws.Cells.Item(42, 5).Value = "  -8.68%  "

$ws.Cells.Item(43, 2).Value = "LidoDAOToken"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$c = $ws.Cells.Item(43, 4)
$c.Value = "'3.38"
$c.Style = "Normal"
$ws.Cells.Item(43, 5).Value = "  +2.93%  "

$ws.Cells.Item(44, 2).Value = "EnergySwap"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$c = $ws.Cells.Item(44, 4)
$c.Value = "'27.35"
$c.Style = "Normal"
$ws.Cells.Item(44, 5).Value = "  +22.72%  "

$ws.Cells.Item(45, 2).Value = "Monero"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$c = $ws.Cells.Item(45, 4)
$c.Value = "'147.89"
$c.Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  +0.47%  "

$ws.Cells.Item(46, 2).Value = "ApeXProtocol"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$c = $ws.Cells.Item(46, 4)
$c.Value = "'3.16"
$c.Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +22.30%  "

$c = $ws.Cells.Item(47, 4)
$c.Value = "'2.09"
$c.Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +3.85%  "

$ws.Cells.Item(48, 2).Value = "NEARProtocol"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$c = $ws.Cells.Item(48, 4)
$c.Value = "'4.35"
$c.Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  -2.65%  "

$ws.Cells.Item(49, 2).Value = "Stacks"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$c = $ws.Cells.Item(49, 4)
$c.Value = "'2.86"
$c.Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  -6.67%  "

$c = $ws.Cells.Item(50, 4)
$c.Value = "'2.60"
$c.Style = "Normal"
$ws.Cells.Item(50, 5).Value = "  -6.56%  "

$c = $ws.Cells.Item(51, 4)
$c.Value = "'0.300"
$c.Style = "Normal"
$ws.Cells.Item(51, 5).Value = "  -3.92%  "

Write-Output "Applied 114 cell changes"